$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# The "名称" (name) field's value was a leftover template placeholder
# ("sheet_productinfo_roomname"); update it to the correct value used in
# this test case.
$ws.Range("B2").Value = "user_info_namenumber"

# Reflect the cell selection left behind by the edit.
$ws.Range("B4").Select()
